$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25-122 down to 26-123.
# Excel's row Insert shifts cell values, formats and styles down with the rows,
# so rows 26-123 end up holding exactly what rows 25-122 held before the insert.
$ws.Rows("25:25").Insert()

# Populate the newly inserted (blank) row 25 with the new record's data.
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C25").Value = "Los Lagos"
$ws.Range("D25").Value = 44701
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 100112052
$ws.Range("G25").Value = "Albahaca"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 7000
$ws.Range("N25").Value = "$/docena de matas"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 1167
$ws.Range("Q25").Value = 6
$ws.Range("R25").Value = "Hortaliza"
